$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 524-525; this pushes the existing rows
# 524-549 down to become 526-551 (matching the dimension change
# from A1:T549 to A1:T551).
$ws.Rows("524:525").Insert()

# New row 524 data
$ws.Cells.Item(524,1).Value = 5
$ws.Cells.Item(524,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(524,3).Value = "Maule"
$ws.Cells.Item(524,4).Value = 44746
$ws.Cells.Item(524,5).Value = 7
$ws.Cells.Item(524,6).Value = "Fruta"
$ws.Cells.Item(524,7).Value = 100102
$ws.Cells.Item(524,8).Value = "Cítricos"
$ws.Cells.Item(524,9).Value = 100102005
$ws.Cells.Item(524,10).Value = "Naranja"
$ws.Cells.Item(524,11).Value = "Fukumoto"
$ws.Cells.Item(524,12).Value = "Primera"
$ws.Cells.Item(524,13).Value = 360
$ws.Cells.Item(524,14).Value = 7000
$ws.Cells.Item(524,15).Value = 7000
$ws.Cells.Item(524,16).Value = 7000
$ws.Cells.Item(524,17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(524,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(524,19).Value = 467
$ws.Cells.Item(524,20).Value = 15

# New row 525 data
$ws.Cells.Item(525,1).Value = 5
$ws.Cells.Item(525,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(525,3).Value = "Maule"
$ws.Cells.Item(525,4).Value = 44746
$ws.Cells.Item(525,5).Value = 7
$ws.Cells.Item(525,6).Value = "Fruta"
$ws.Cells.Item(525,7).Value = 100102
$ws.Cells.Item(525,8).Value = "Cítricos"
$ws.Cells.Item(525,9).Value = 100102005
$ws.Cells.Item(525,10).Value = "Naranja"
$ws.Cells.Item(525,11).Value = "New Hall"
$ws.Cells.Item(525,12).Value = "Primera"
$ws.Cells.Item(525,13).Value = 560
$ws.Cells.Item(525,14).Value = 6500
$ws.Cells.Item(525,15).Value = 7000
$ws.Cells.Item(525,16).Value = 6821
$ws.Cells.Item(525,17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(525,18).Value = "Región de O'Higgins"
$ws.Cells.Item(525,19).Value = 455
$ws.Cells.Item(525,20).Value = 15
